# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
# A new weekly record is inserted as row 60 (pushing the existing rows
# 60-147 down to 61-148), carrying forward the same market/category
# metadata as the row that used to occupy row 60, but with an updated
# date and adjusted volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60; everything below shifts down one row.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Cells.Item(60, 1).Value  = 1
$ws.Cells.Item(60, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(60, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(60, 4).Value  = 45100
$ws.Cells.Item(60, 5).Value  = 15
$ws.Cells.Item(60, 6).Value  = 100112036
$ws.Cells.Item(60, 7).Value  = "Caigua"
$ws.Cells.Item(60, 8).Value  = "Sin especificar"
$ws.Cells.Item(60, 9).Value  = "Primera"
$ws.Cells.Item(60, 10).Value = 130
$ws.Cells.Item(60, 11).Value = 11000
$ws.Cells.Item(60, 12).Value = 12000
$ws.Cells.Item(60, 13).Value = 11500
$ws.Cells.Item(60, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(60, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(60, 16).Value = 767
$ws.Cells.Item(60, 17).Value = 15
$ws.Cells.Item(60, 18).Value = "Hortaliza"

# Match the date-number-format style used by the rest of column D.
$ws.Cells.Item(60, 4).NumberFormat = $ws.Cells.Item(61, 4).NumberFormat
